$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "0.1.2"
$meta.Range("B5").Value = "CodeSystem - Blood Group (Rh) - NMDP"
$meta.Range("B8").Value = "2025-04-16T10:37:17-05:00"

# "Count" (B22) must stay text even though "2" looks numeric (it was a
# shared string before the edit too). Force text entry, then restore the
# original cell formatting (copied from the cell above) so only the
# value itself changes.
$meta.Range("B22").NumberFormat = "@"
$meta.Range("B22").Value = "2"
$meta.Range("B21").Copy()
$meta.Range("B22").PasteSpecial(-4122)

# --- Concepts sheet updates ---
$concepts = $wb.Worksheets.Item("Concepts")

# Update remaining rows to the new Rh+/Rh- codes
$concepts.Range("B2").Value = "Rh+"
$concepts.Range("C2").Value = "Positive"
$concepts.Range("B3").Value = "Rh-"
$concepts.Range("C3").Value = "Negative"

# Remove the now-obsolete rows (I/Indeterminant, D/nil/no data, etc.)
$concepts.Rows.Item(4).Delete()
$concepts.Rows.Item(4).Delete()
$concepts.Rows.Item(4).Delete()
